$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 5.3
$ws.Range("P2").Value = 1.89
$ws.Range("Q3").Value = 2.16
$ws.Range("G4").Value = 1.48
$ws.Range("J4").Value = 4.4
$ws.Range("F5").Value = 1.99
$ws.Range("G5").Value = 2.24
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 4.1
$ws.Range("J6").Value = 4.1
$ws.Range("Q7").Value = 1.01
$ws.Range("AE8").Value = 18
$ws.Range("AG8").Value = 18
$ws.Range("AH8").Value = 19.5
$ws.Range("AN8").Value = 65
$ws.Range("F8").Value = 4.3
$ws.Range("G8").Value = 4.6
$ws.Range("H8").Value = 1.83
$ws.Range("I8").Value = 1.86
$ws.Range("Q8").Value = 1.78
$ws.Range("T8").Value = 1.77
$ws.Range("U8").Value = 2.18
$ws.Range("Z8").Value = 11.5
$ws.Range("T10").Value = 1.88
$ws.Range("AC11").Value = 13.5
$ws.Range("AD11").Value = 42
$ws.Range("AE11").Value = 240
$ws.Range("AH11").Value = 34
$ws.Range("AJ11").Value = 10.5
$ws.Range("AL11").Value = 42
$ws.Range("AM11").Value = 240
$ws.Range("F11").Value = 1.35
$ws.Range("G11").Value = 1.36
$ws.Range("H11").Value = 10.5
$ws.Range("J11").Value = 5.8
$ws.Range("K11").Value = 6
$ws.Range("N11").Value = 4.6
$ws.Range("O11").Value = 1.25
$ws.Range("P11").Value = 2.26
$ws.Range("S11").Value = 2.88
$ws.Range("Y11").Value = 34
$ws.Range("AF12").Value = 40
$ws.Range("AI12").Value = 40
$ws.Range("AN12").Value = 110
$ws.Range("G12").Value = 5.2
$ws.Range("K12").Value = 4
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 3.6
$ws.Range("V13").Value = 1.38
$ws.Range("H14").Value = 2.58
$ws.Range("P14").Value = 1.9
$ws.Range("Q14").Value = 1.64
$ws.Range("AO15").Value = 9
$ws.Range("I15").Value = 1.79
$ws.Range("P15").Value = 2.2
$ws.Range("Q15").Value = 1.69
$ws.Range("S15").Value = 2.68
$ws.Range("T15").Value = 1.72
$ws.Range("U15").Value = 2.12
$ws.Range("L16").Value = 1.26
$ws.Range("U16").Value = 1.79
$ws.Range("AB17").Value = 1000
$ws.Range("AF17").Value = 1000
$ws.Range("AG17").Value = 1000
$ws.Range("AJ17").Value = 1000
$ws.Range("AK17").Value = 1000
$ws.Range("AN17").Value = 1000
$ws.Range("F17").Value = 1.1
$ws.Range("H17").Value = 1.06
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 1.01
$ws.Range("N17").Value = 7.6
$ws.Range("O17").Value = 1.08
$ws.Range("P17").Value = 3.6
$ws.Range("Q17").Value = 1.23
$ws.Range("R17").Value = 1.44
$ws.Range("S17").Value = 1.23
$ws.Range("T17").Value = 1.01
$ws.Range("U17").Value = 1.01
$ws.Range("W17").Value = 5.4
$ws.Range("F18").Value = 2.3
$ws.Range("G18").Value = 2.88
$ws.Range("H18").Value = 2.82
$ws.Range("I18").Value = 3.8
$ws.Range("J18").Value = 3.1
$ws.Range("N18").Value = 3.3
$ws.Range("O18").Value = 1.24
$ws.Range("P18").Value = 1.79
$ws.Range("Q18").Value = 1.87
$ws.Range("R18").Value = 1.26
$ws.Range("S18").Value = 2.86
$ws.Range("V18").Value = 1.36
$ws.Range("W18").Value = 1.53
$ws.Range("AI19").Value = 110
$ws.Range("AN19").Value = 5.2
$ws.Range("I19").Value = 9.199999999999999
$ws.Range("AH20").Value = 20
$ws.Range("AK20").Value = 550
$ws.Range("AL20").Value = 75
$ws.Range("F20").Value = 4.7
$ws.Range("G20").Value = 4.8
$ws.Range("H20").Value = 1.82
$ws.Range("I20").Value = 1.83
$ws.Range("P20").Value = 2.04
$ws.Range("R20").Value = 1.41
$ws.Range("V20").Value = 2.2
$ws.Range("W20").Value = 1.26
$ws.Range("F21").Value = 2.04
$ws.Range("G21").Value = 2.16
$ws.Range("F22").Value = 2.36
$ws.Range("Q22").Value = 1.86
$ws.Range("I23").Value = 1.85
$ws.Range("P23").Value = 1.86
$ws.Range("Q23").Value = 1.99
$ws.Range("H25").Value = 2.9
$ws.Range("I25").Value = 3.1
$ws.Range("K25").Value = 3.5
$ws.Range("P25").Value = 1.65
$ws.Range("Q25").Value = 2.26
